$d = $word.ActiveDocument

$pairs = @(
    @("902×8=", "900×7="),
    @("683×8=", "566×7="),
    @("560×2=", "980×7="),
    @("965×6=", "519×4="),
    @("442×9=", "528×6="),
    @("246×5=", "300×7="),
    @("511×2=", "389×6="),
    @("531×2=", "123×8="),
    @("437×8=", "350×2="),
    @("982×3=", "257×8="),
    @("713×5=", "748×7="),
    @("122×2=", "577×4="),
    @("823×4=", "755×4="),
    @("526×9=", "183×5="),
    @("107×7=", "395×7="),
    @("639×8=", "237×2="),
    @("708×2=", "171×6="),
    @("743×4=", "812×8="),
    @("873×4=", "474×8="),
    @("878×4=", "136×7="),
    @("423×2=", "710×9="),
    @("530×9=", "852×6="),
    @("775×7=", "395×2="),
    @("435×5=", "151×5="),
    @("973×7=", "297×7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
